$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 436.66666
$ws.Range("I5").Value = 436.66666
$ws.Range("K5").Value = 436.66666
$ws.Range("M5").Value = -321.66666
$ws.Range("H28").Value = 71429224
$ws.Range("I28").Value = 90909620
$ws.Range("K28").Value = 90909620
$ws.Range("M28").Value = -90909135
$ws.Range("H41").Value = 412.91666
$ws.Range("I41").Value = 407
$ws.Range("K41").Value = 407
$ws.Range("M41").Value = 33
$ws.Range("H53").Value = 381.30435
$ws.Range("I53").Value = 515.2308
$ws.Range("J53").Value = 207.2
$ws.Range("K53").Value = 515.2308
$ws.Range("L53").Value = 207.2
$ws.Range("M53").Value = 121.7692
$ws.Range("N53").Value = -1481.2
$ws.Range("H55").Value = 150.28572
$ws.Range("I55").Value = 89.818184
$ws.Range("J55").Value = 372
$ws.Range("K55").Value = 89.818184
$ws.Range("L55").Value = 372
$ws.Range("M55").Value = 124.181816
$ws.Range("N55").Value = -800
$ws.Range("H57").Value = 134335.4
$ws.Range("J57").Value = 134335.4
$ws.Range("L57").Value = 403006.2
$ws.Range("N57").Value = -404004.2
$ws.Range("H69").Value = 18607.084
$ws.Range("I69").Value = 7715.5
$ws.Range("J69").Value = 29498.666
$ws.Range("K69").Value = 23146.5
$ws.Range("L69").Value = 88495.99800000001
$ws.Range("M69").Value = -22272.5
$ws.Range("N69").Value = -90243.99800000001
$ws.Range("H72").Value = 18607.084
$ws.Range("I72").Value = 7715.5
$ws.Range("J72").Value = 29498.666
$ws.Range("K72").Value = 69439.5
$ws.Range("L72").Value = 265487.994
$ws.Range("M72").Value = -65071.5
$ws.Range("N72").Value = -274223.994
$ws.Range("H107").Value = 41668544
$ws.Range("I107").Value = 62502292
$ws.Range("J107").Value = 1049.25
$ws.Range("K107").Value = 62502292
$ws.Range("L107").Value = 1049.25
$ws.Range("M107").Value = -62500372
$ws.Range("N107").Value = -4889.25
$ws.Range("H111").Value = 926
$ws.Range("I111").Value = 718.7778
$ws.Range("J111").Value = 1299
$ws.Range("K111").Value = 2156.3334
$ws.Range("L111").Value = 3897
$ws.Range("M111").Value = 910.6666
$ws.Range("N111").Value = -10031
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1585.8966
$ws.Range("I2").Value = 1528.64
$ws.Range("J2").Value = 1943.75
$ws.Range("K2").Value = 1528.64
$ws.Range("L2").Value = 1943.75
$ws.Range("M2").Value = -1415.64
$ws.Range("N2").Value = -2169.75
$ws.Range("H4").Value = 505.2
$ws.Range("I4").Value = 202.22223
$ws.Range("K4").Value = 202.22223
$ws.Range("M4").Value = -86.22223
$ws.Range("H5").Value = 193.8
$ws.Range("I5").Value = 180
$ws.Range("J5").Value = 249
$ws.Range("K5").Value = 180
$ws.Range("L5").Value = 249
$ws.Range("M5").Value = -68
$ws.Range("N5").Value = -473
$ws.Range("H110").Value = 3306.3215
$ws.Range("I110").Value = 3243.5925
$ws.Range("K110").Value = 3243.5925
$ws.Range("M110").Value = -1198.5925
$ws.Range("H116").Value = 1585.8966
$ws.Range("I116").Value = 1528.64
$ws.Range("J116").Value = 1943.75
$ws.Range("K116").Value = 1528.64
$ws.Range("L116").Value = 1943.75
$ws.Range("M116").Value = 765.3599999999999
$ws.Range("N116").Value = -6531.75
$ws.Range("H132").Value = 2781.1
$ws.Range("I132").Value = 2607.2856
$ws.Range("J132").Value = 3186.6667
$ws.Range("K132").Value = 7821.8568
$ws.Range("L132").Value = 9560.000100000001
$ws.Range("M132").Value = -5291.8568
$ws.Range("N132").Value = -14620.0001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1585.8966
$ws.Range("I3").Value = 1528.64
$ws.Range("J3").Value = 1943.75
$ws.Range("K3").Value = 1528.64
$ws.Range("L3").Value = 1943.75
$ws.Range("M3").Value = -1414.64
$ws.Range("N3").Value = -2171.75
$ws.Range("H4").Value = 193.8
$ws.Range("I4").Value = 180
$ws.Range("J4").Value = 249
$ws.Range("K4").Value = 180
$ws.Range("L4").Value = 249
$ws.Range("M4").Value = -65
$ws.Range("N4").Value = -479
$ws.Range("H80").Value = 2882.2083
$ws.Range("J80").Value = 881.2
$ws.Range("L80").Value = 881.2
$ws.Range("N80").Value = -2877.2
$ws.Range("H83").Value = 2882.2083
$ws.Range("J83").Value = 881.2
$ws.Range("L83").Value = 4406
$ws.Range("N83").Value = -14390
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 68.57143000000001
$ws.Range("I7").Value = 62.6
$ws.Range("J7").Value = 83.5
$ws.Range("K7").Value = 62.6
$ws.Range("L7").Value = 83.5
$ws.Range("M7").Value = 50.4
$ws.Range("N7").Value = -309.5
$ws.Range("H9").Value = 99932
$ws.Range("J9").Value = 99932
$ws.Range("L9").Value = 99932
$ws.Range("N9").Value = -100268
$ws.Range("H22").Value = 957.69696
$ws.Range("I22").Value = 697.6
$ws.Range("K22").Value = 697.6
$ws.Range("M22").Value = -347.6
$ws.Range("H86").Value = 3341693.5
$ws.Range("I86").Value = 6066315.5
$ws.Range("J86").Value = 11599.777
$ws.Range("K86").Value = 6066315.5
$ws.Range("L86").Value = 11599.777
$ws.Range("M86").Value = -6065192.5
$ws.Range("N86").Value = -13845.777
$ws.Range("H89").Value = 3341693.5
$ws.Range("I89").Value = 6066315.5
$ws.Range("J89").Value = 11599.777
$ws.Range("K89").Value = 30331577.5
$ws.Range("L89").Value = 57998.885
$ws.Range("M89").Value = -30325961.5
$ws.Range("N89").Value = -69230.88500000001
$ws.Range("H94").Value = 2708.5
$ws.Range("I94").Value = 1503.6666
$ws.Range("K94").Value = 1503.6666
$ws.Range("M94").Value = -1052.6666
$ws.Range("H99").Value = 7327.4688
$ws.Range("I99").Value = 5867.8667
$ws.Range("K99").Value = 5867.8667
$ws.Range("M99").Value = -4369.8667
$ws.Range("H126").Value = 7327.4688
$ws.Range("I126").Value = 5867.8667
$ws.Range("K126").Value = 17603.6001
$ws.Range("M126").Value = -15133.6001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1288.2858
$ws.Range("I113").Value = 203
$ws.Range("J113").Value = 1469.1666
$ws.Range("K113").Value = 609
$ws.Range("L113").Value = 4407.4998
$ws.Range("N113").Value = -8747.4998
$ws.Range("M113").Value = 1561
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 259259420
$ws.Range("I2").Value = 259259420
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 259259420
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -259259307
$ws.Range("N2").ClearContents()
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H135").Value = 191998.14
$ws.Range("J135").Value = 191998.14
$ws.Range("L135").Value = 191998.14
$ws.Range("N135").Value = -202138.14
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3833.0967
$ws.Range("I7").Value = 3549.476
$ws.Range("J7").Value = 4428.7
$ws.Range("K7").Value = 3549.476
$ws.Range("L7").Value = 4428.7
$ws.Range("M7").Value = -3437.476
$ws.Range("N7").Value = -4652.7
$ws.Range("H82").Value = 658.64514
$ws.Range("J82").Value = 1882.3334
$ws.Range("L82").Value = 1882.3334
$ws.Range("N82").Value = -2604.3334
$ws.Range("H85").Value = 658.64514
$ws.Range("J85").Value = 1882.3334
$ws.Range("L85").Value = 1882.3334
$ws.Range("N85").Value = -4378.3334
$ws.Range("H93").Value = 1953.4667
$ws.Range("I93").Value = 1976.037
$ws.Range("J93").Value = 1750.3334
$ws.Range("K93").Value = 1976.037
$ws.Range("L93").Value = 1750.3334
$ws.Range("M93").Value = -728.037
$ws.Range("N93").Value = -4246.3334
$ws.Range("H126").Value = 3833.0967
$ws.Range("I126").Value = 3549.476
$ws.Range("J126").Value = 4428.7
$ws.Range("K126").Value = 10648.428
$ws.Range("L126").Value = 13286.1
$ws.Range("M126").Value = -8178.428
$ws.Range("N126").Value = -18226.1
$ws.Range("H137").Value = 99761.336
$ws.Range("J137").Value = 129947
$ws.Range("L137").Value = 129947
$ws.Range("N137").Value = -140147
$ws.Range("H139").Value = 99995.5
$ws.Range("J139").Value = 99995.5
$ws.Range("L139").Value = 99995.5
$ws.Range("N139").Value = -110275.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 8221.5
$ws.Range("J19").Value = 9443
$ws.Range("L19").Value = 9443
$ws.Range("N19").Value = -9791
$ws.Range("H62").Value = 9098.299999999999
$ws.Range("I62").Value = 8898.833000000001
$ws.Range("J62").Value = 9397.5
$ws.Range("K62").Value = 8898.833000000001
$ws.Range("L62").Value = 9397.5
$ws.Range("M62").Value = -8274.833000000001
$ws.Range("N62").Value = -10645.5
$ws.Range("H65").Value = 9098.299999999999
$ws.Range("I65").Value = 8898.833000000001
$ws.Range("J65").Value = 9397.5
$ws.Range("K65").Value = 44494.165
$ws.Range("L65").Value = 46987.5
$ws.Range("M65").Value = -41374.165
$ws.Range("N65").Value = -53227.5
$ws.Range("H122").Value = 4301.074
$ws.Range("I122").Value = 4702.263
$ws.Range("J122").Value = 3348.25
$ws.Range("K122").Value = 14106.789
$ws.Range("L122").Value = 10044.75
$ws.Range("M122").Value = -11656.789
$ws.Range("N122").Value = -14944.75
$ws.Range("H126").Value = 5280.7144
$ws.Range("I126").Value = 4669.0415
$ws.Range("J126").Value = 8950.75
$ws.Range("K126").Value = 14007.1245
$ws.Range("L126").Value = 26852.25
$ws.Range("M126").Value = -11537.1245
$ws.Range("N126").Value = -31792.25
$ws.Range("H132").Value = 3498.5
$ws.Range("I132").Value = 3499.5
$ws.Range("K132").Value = 10498.5
$ws.Range("M132").Value = -7968.5
